$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 4")

# Remove the "enemy info" task row (old row 15) -- rows below shift up by one.
$ws.Rows("15:15").Delete()

# Program Display / weapon-effect task: status + notes updated for Michael's contribution.
$ws.Range("F5").Value = "done"
$ws.Range("G5").Value = "different weapons are implemented"

# Enemies / spawn task: status + notes updated.
$ws.Range("F10").Value = "done"
$ws.Range("G10").Value = "the player is damaged when hit by a enemy shot"

# Level / difficulty-settings task: status + notes updated.
$ws.Range("F12").Value = "done"
$ws.Range("G12").Value = "weapons do differing amounts of damage to the ships"

# Level / "read what each weapon does" task: status + notes updated.
$ws.Range("F14").Value = "done"
$ws.Range("G14").Value = "all items have over behavior for them"

# Audio task notes updated (row shifted from 20 to 19 after the deletion above).
$ws.Range("G19").Value = "audio does play but the volume is not controllable, not mutes"

# Restore the view state recorded after these edits.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("F14").Select()
